$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.894.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3892'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3845'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.02'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.355'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08469'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.174'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.904'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001305'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.654.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07010'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.919'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.896.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.503'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.044'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.416'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '139.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.827'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.492'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.836.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.036'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08073'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02966'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.15%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.80%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.681'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2690'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09132'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7533'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.419'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6961'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.467'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.078'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08271'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.237'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
